# "CHanges of 30th Jan"
# Refresh the Service IDs (col C) and Fail Log text (col F) for the rows that
# were re-run, same pattern as the existing sheet: plain text/shared-string
# values, no formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    # Writes $text into $cellAddr as a genuine text value (not auto-coerced
    # to a number by Excel) while leaving the cell's visible style untouched.
    param($cellAddr, $text)
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Column C - Service ID values
Set-TextValue "C2" "10262084"
Set-TextValue "C3" "10262092"
Set-TextValue "C4" "10262104"
Set-TextValue "C5" "10262110"
Set-TextValue "C11" "10262213"
Set-TextValue "C12" "10262208"
Set-TextValue "C13" "10262150"
Set-TextValue "C14" "10262211"
Set-TextValue "C24" "136839363"

# Column F - Fail Log values
$seleniumNoSuchElement = @"
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome-headless-shell=121.0.6167.85)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.69', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '19.0.1'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 121.0.6167.85, chrome: {chromedriverVersion: 121.0.6167.85 (3f98d690ad7e..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:64237}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: d4804c12644307a50394ee1d762ff46a
*** Element info: {Using=id, value=lblServiceID}
"@

$ws.Range("F11").Value = $seleniumNoSuchElement
$ws.Range("F12").Value = $seleniumNoSuchElement
$ws.Range("F14").Value = 'Cannot invoke "org.openqa.selenium.WebElement.isDisplayed()" because "element" is null'

# The Selenium dump is multi-line and column F wraps text, so writing it
# bumps the row to a custom height; auto-fit puts it back to the sheet
# default (matches the unedited rows, and the diff shows no row changes).
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(14).AutoFit()
